# "changes in mobile view" - add two new tenant rows (row 3 and row 4) to the
# Tenants sheet, mirroring the existing header/row-2 layout.
#
# Columns: A ID | B Full Name | C Company | D Contact | E Citizen Number |
#          F Address | G Property ID | H Section | I Contract Years |
#          J Start Date AD | K Start Date BS | L End Date AD | M End Date BS |
#          N Amount | O Amount Type | P Increment Percent |
#          Q Increment Interval | R Status

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Force literal text storage (no auto date/number detection) using the
    # classic quote-prefix trick, then drop the resulting "quote prefix"
    # cell style back to Normal so no visible formatting change sticks.
    $ws.Range($range).Value = "'" + $text
    $ws.Range($range).Style = "Normal"
}

# ---- Row 3 : Sudip Maharjan ----
Set-TextValue "A3" "1756651364942"
$ws.Range("B3").Value = "Sudip Maharjan"
Set-TextValue "C3" ""
Set-TextValue "D3" "9865072119"
Set-TextValue "E3" ""
$ws.Range("F3").Value = "Ravi Bhawan"
Set-TextValue "G3" "1755340190541"
$ws.Range("I3").Value = 5
Set-TextValue "J3" "2025-08-01"
Set-TextValue "K3" "2025-08-01"
Set-TextValue "L3" "2030-07-31"
Set-TextValue "M3" "2030-07-31"
$ws.Range("N3").Value = 50000
$ws.Range("O3").Value = "month"
$ws.Range("P3").Value = 10
$ws.Range("Q3").Value = 2
$ws.Range("R3").Value = "Active"

# ---- Row 4 : XYZ / ABC ----
Set-TextValue "A4" "1756651403346"
$ws.Range("B4").Value = "XYZ"
$ws.Range("C4").Value = "ABC"
Set-TextValue "D4" "98562"
Set-TextValue "E4" "2"
$ws.Range("F4").Value = "Kulshwor"
Set-TextValue "G4" "1755500229573"
$ws.Range("H4").Value = "1st Floor (201)"
$ws.Range("I4").Value = 10
Set-TextValue "J4" "2025-09-01"
Set-TextValue "K4" "2025-09-01"
Set-TextValue "L4" "2035-08-31"
Set-TextValue "M4" "2035-08-31"
$ws.Range("N4").Value = 24999
$ws.Range("O4").Value = "month"
$ws.Range("P4").Value = 10
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = "Active"
